$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" field text from 08/17/2020 to
#    08/18/2020 on the slide master and on every slide layout's date
#    placeholder.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Type -eq 14 -and $shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                if ($shp.TextFrame.TextRange.Text -eq "08/17/2020") {
                    $shp.TextFrame.TextRange.Text = "08/18/2020"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder($master.Shapes)

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder($layout.Shapes)
}

# ---------------------------------------------------------------------------
# 2) Add the author's name + student number as two new text boxes on the
#    title slide (slide 1), in the top-right corner.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

$nameBox = $slide1.Shapes.AddTextbox(1, 646.507716535433, 40.97496062992126, 148.3548031496063, 29.081259842519685)
$nameBox.Name = "Tekstfelt 3"
$nameBox.Fill.Visible = $false
$nameBox.TextFrame.WordWrap = $false
$nameBox.TextFrame.AutoSize = 1
$nameBox.TextFrame.TextRange.Text = "Andreas Blaabjerg"
$nameBox.TextFrame.TextRange.LanguageID = "da-DK"

$idBox = $slide1.Shapes.AddTextbox(1, 720.6851181102362, 72.29228346456694, 97.46763779527559, 29.081259842519685)
$idBox.Name = "Tekstfelt 4"
$idBox.Fill.Visible = $false
$idBox.TextFrame.WordWrap = $false
$idBox.TextFrame.AutoSize = 1
$idBox.TextFrame.TextRange.Text = "201510924"
$idBox.TextFrame.TextRange.LanguageID = "da-DK"

Write-Host "done"
